$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new "Random" sample row appended under the existing header ---

# A2: the date/time value, formatted with the built-in m/d/yy h:mm number format
# (this mints a new cellXf with numFmtId 22, applyNumberFormat="1")
$ws.Range("A2").Value = 42609.647106481483
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"

# B2: plain integer
$ws.Range("B2").Value = 39

# C2:M2: zeros
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0

# N2: new shared string "Random"
$ws.Range("N2").Value = "Random"

# A1 picks up the same date style as A2 (column-header cell re-styled too)
$ws.Range("A1").NumberFormat = "m/d/yy h:mm"

# Column A widens to fit the new date/time column contents
$ws.Columns.Item(1).ColumnWidth = 14
